# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetExpo = $wb.Worksheets.Item("展览")
$sheetAll  = $wb.Worksheets.Item("全部类型")

# Row -> (old, new) mapping for the "展览" sheet
$sheetExpo.Range("F7").Value  = 2706
$sheetExpo.Range("F9").Value  = 574
$sheetExpo.Range("F12").Value = 10287
$sheetExpo.Range("F14").Value = 272
$sheetExpo.Range("F16").Value = 645
$sheetExpo.Range("F17").Value = 11859
$sheetExpo.Range("F18").Value = 12251
$sheetExpo.Range("F20").Value = 103
$sheetExpo.Range("F22").Value = 27

# Row -> (old, new) mapping for the "全部类型" sheet
$sheetAll.Range("F7").Value  = 2706
$sheetAll.Range("F10").Value = 574
$sheetAll.Range("F13").Value = 10287
$sheetAll.Range("F15").Value = 272
$sheetAll.Range("F17").Value = 645
$sheetAll.Range("F18").Value = 11859
$sheetAll.Range("F19").Value = 12251
$sheetAll.Range("F21").Value = 103
$sheetAll.Range("F23").Value = 27
